$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Maps 2 RGB")

# --- Re-run results for V4.06.5 .. V4.06.8 (rows 25-28): new numbers, and the
#     red "flagged" font is replaced with the normal numeric style already
#     used elsewhere in the workbook (e.g. Sheet2!F10) for this number format.
$plainNumberFormat = $wb.Worksheets.Item("Sheet2").Range("F10")
$plainNumberFormat.Copy()
$ws.Range("B25:I28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B25").Value = 14.998900000000001
$ws.Range("C25").Value = 0.73729999999999996
$ws.Range("D25").Value = 13.646100000000001
$ws.Range("E25").Value = 0.79569999999999996
$ws.Range("F25").Value = 17.9329
$ws.Range("G25").Value = 0.84189999999999998
$ws.Range("H25").Value = 20.722300000000001
$ws.Range("I25").Value = 0.81699999999999995

$ws.Range("B26").Value = 9.1883999999999997
$ws.Range("C26").Value = 0.65439999999999998
$ws.Range("D26").Value = 13.919
$ws.Range("E26").Value = 0.79610000000000003
$ws.Range("F26").Value = 17.896899999999999
$ws.Range("G26").Value = 0.84860000000000002
$ws.Range("H26").Value = 18.140499999999999
$ws.Range("I26").Value = 0.82150000000000001

$ws.Range("B27").Value = 7.6097999999999999
$ws.Range("C27").Value = 0.62519999999999998
$ws.Range("D27").Value = 14.295199999999999
$ws.Range("E27").Value = 0.81679999999999997
$ws.Range("F27").Value = 18.012899999999998
$ws.Range("G27").Value = 0.83309999999999995
$ws.Range("H27").Value = 19.110399999999998
$ws.Range("I27").Value = 0.81699999999999995

$ws.Range("B28").Value = 7.6630000000000003
$ws.Range("C28").Value = 0.62439999999999996
$ws.Range("D28").Value = 14.385300000000001
$ws.Range("E28").Value = 0.82469999999999999
$ws.Range("F28").Value = 19.760200000000001
$ws.Range("G28").Value = 0.85970000000000002
$ws.Range("H28").Value = 19.948699999999999
$ws.Range("I28").Value = 0.8327

# --- V4.07.5 (row 29): the run hasn't been done yet, so its results are
#     cleared out and the row is given the same "pending" centered style used
#     for the group headers (e.g. Sheet1!A1), instead of the red flagged one.
$pendingStyle = $wb.Worksheets.Item("Sheet1").Range("A1")
$pendingStyle.Copy()
$ws.Range("B29:I29").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B29:I29").ClearContents()

# --- Scroll position / selection, matching where the author left off editing.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D30").Select() | Out-Null
